# "from home PC main:Branch"
#
# The "Merge Result" sheet gains the extra columns that come from merging in
# the "State" column (plus the trailing placeholder Column1..Column12 fields
# that the merge produced) from the "Sheet2" table. Only the two rows that
# had a matching State value end up populated in column D; the rest of the
# new columns (E:P) stay empty except for their header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Merge Result")

# New header row: column D onward, starting right after the existing
# Name / ID / Place headers in A1:C1.
$headers = @(
    "State",
    "Column1", "Column2", "Column3", "Column4", "Column5", "Column6",
    "Column7", "Column8", "Column9", "Column10", "Column11", "Column12"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, 4 + $i).Value = $headers[$i]
}

# Only rows 4 and 5 (the duplicated pavan/sandeep rows) pick up a matched
# State value from the merge.
$ws.Cells.Item(4, 4).Value = "AP"
$ws.Cells.Item(5, 4).Value = "TS"
